# Update "想去人数" (F column) counts across the 展览, 演出 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (index 1) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 97
$ws1.Range("F3").Value  = 171
$ws1.Range("F4").Value  = 406
$ws1.Range("F5").Value  = 184
$ws1.Range("F6").Value  = 128
$ws1.Range("F7").Value  = 1078
$ws1.Range("F8").Value  = 359
$ws1.Range("F9").Value  = 183
$ws1.Range("F11").Value = 120
$ws1.Range("F12").Value = 366
$ws1.Range("F13").Value = 364
$ws1.Range("F14").Value = 775
$ws1.Range("F15").Value = 148
$ws1.Range("F16").Value = 714
$ws1.Range("F17").Value = 269
$ws1.Range("F19").Value = 986
$ws1.Range("F20").Value = 444
$ws1.Range("F21").Value = 254
$ws1.Range("F23").Value = 371
$ws1.Range("F25").Value = 38

# --- Sheet 2: 演出 (index 2) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F11").Value = 147

# --- Sheet 4: 全部类型 (index 4) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 97
$ws4.Range("F5").Value  = 171
$ws4.Range("F6").Value  = 406
$ws4.Range("F7").Value  = 184
$ws4.Range("F8").Value  = 128
$ws4.Range("F9").Value  = 1078
$ws4.Range("F10").Value = 359
$ws4.Range("F11").Value = 183
$ws4.Range("F15").Value = 120
$ws4.Range("F17").Value = 366
$ws4.Range("F20").Value = 364
$ws4.Range("F21").Value = 775
$ws4.Range("F22").Value = 148
$ws4.Range("F23").Value = 714
$ws4.Range("F24").Value = 269
$ws4.Range("F26").Value = 986
$ws4.Range("F27").Value = 444
$ws4.Range("F30").Value = 254
$ws4.Range("F32").Value = 371
$ws4.Range("F34").Value = 147
$ws4.Range("F36").Value = 38
